$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("J62").Value = 0.148
$ws.Range("G65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("G68").Value = 0.02500000000000002
$ws.Range("K68").Value = -0.04699999999999999
$ws.Range("J81").Value = 0.123
$ws.Range("J84").Value = 0.135
$ws.Range("J88").Value = 0.02399999999999997
$ws.Range("K88").Value = -0.03599999999999998
$ws.Range("E89").Value = -0.03600000000000003
$ws.Range("I89").Value = 0.01600000000000001
$ws.Range("J93").Value = 0.04100000000000004
$ws.Range("G97").Value = 0.09300000000000008
$ws.Range("L99").Value = -0.01299999999999996
$ws.Range("J100").Value = 0.04799999999999999
